# Append the latest EUR->ARS quote as a new row at the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 97

# Column A holds a date-looking string ("2025-10-24"). Excel's COM layer
# auto-converts such text to a date serial on assignment unless the cell is
# already formatted as Text, so force that first (matches how the sheet's
# existing "Fecha" values are stored as literal strings).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-10-24"

$ws.Cells.Item($row, 2).Value = "15:28:30"
$ws.Cells.Item($row, 3).Value = "1.00 EUR = 1,828.9067"
